# Update cryptos list with latest scraped prices / volume percentages.
# D-column price strings are stored as TEXT (not numbers) in the source
# data, so a leading apostrophe is used to force Excel to keep them as
# text instead of auto-converting numeric-looking strings (e.g. "211.42")
# into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "26.667.65"; E = "  +0.01%  " },
    @{ Row = 3;  D = "1.598.02";  E = "  +0.03%  " },
    @{ Row = 4;  E = "  +0.16%  " },
    @{ Row = 5;  D = "211.42";    E = "  +0.05%  " },
    @{ Row = 6;  D = "0.515";     E = "  +0.76%  " },
    @{ Row = 9;  E = "  +0.79%  " },
    @{ Row = 10; E = "  -0.74%  " },
    @{ Row = 11; D = "0.0840";    E = "  +0.20%  " },
    @{ Row = 12; D = "1.822.01";  E = "  +0.04%  " },
    @{ Row = 13; D = "1.605.86";  E = "  +0.42%  " },
    @{ Row = 14; E = "  -0.12%  " },
    @{ Row = 15; D = "0.524";     E = "  +0.47%  " },
    @{ Row = 16; D = "65.20" },
    @{ Row = 17; D = "26.655.45"; E = "  +0.00%  " },
    @{ Row = 19; D = "209.82";    E = "  +0.01%  " },
    @{ Row = 20; E = "  +0.19%  " },
    @{ Row = 21; D = "7.02";      E = "  +3.75%  " },
    @{ Row = 22; E = "  +0.77%  " },
    @{ Row = 23; E = "  +1.76%  " },
    @{ Row = 24; D = "8.98";      E = "  +0.75%  " },
    @{ Row = 25; E = "  -1.50%  " },
    @{ Row = 26; E = "  +0.11%  " },
    @{ Row = 27; D = "7.12";      E = "  -0.77%  " },
    @{ Row = 28; E = "  -0.39%  " },
    @{ Row = 29; E = "  -0.09%  " },
    @{ Row = 30; E = "  +2.43%  " },
    @{ Row = 31; E = "  +0.50%  " },
    @{ Row = 32; D = "3.25";      E = "  +0.89%  " },
    @{ Row = 33; E = "  +1.63%  " },
    @{ Row = 34; D = "1.287.46" },
    @{ Row = 35; D = "0.620";     E = "  -6.64%  " },
    @{ Row = 36; E = "  +0.45%  " },
    @{ Row = 37; E = "  +0.79%  " },
    @{ Row = 38; E = "  -0.71%  " },
    @{ Row = 39; D = "0.835";     E = "  -0.96%  " },
    @{ Row = 40; E = "  +20.23%  " },
    @{ Row = 41; E = "  +2.32%  " },
    @{ Row = 42; E = "  -0.25%  " },
    @{ Row = 43; D = "0.786";     E = "  -0.33%  " },
    @{ Row = 44; D = "63.57";     E = "  -0.29%  " },
    @{ Row = 45; D = "1.735.34";  E = "  +0.08%  " },
    @{ Row = 46; D = "90.79";     E = "  +0.65%  " },
    @{ Row = 47; E = "  -3.51%  " },
    @{ Row = 48; E = "  -0.03%  " },
    @{ Row = 49; E = "  +1.14%  " },
    @{ Row = 50; D = "0.0509";    E = "  +0.87%  " },
    @{ Row = 51; E = "  +0.15%  " }
)

$quote = "'"

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        # Leading apostrophe forces text storage so numeric-looking
        # strings aren't coerced into Number cells.
        $ws.Cells.Item($r, 4).Value = $quote + $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
